$d = $word.ActiveDocument

# --- Move the "_GoBack" bookmark from the heading paragraph down into the
# --- "section and team ID" sentence, right after "001_03" (this is what
# --- Word does automatically to mark the most-recent edit location).
$oldGoBack = $d.Bookmarks.Item("_GoBack")
$oldGoBack.Delete()

$findRange = $d.Content
[void]$findRange.Find.Execute("e.g., 001_03", $false, $false, $false, $false, $false, `
                               $true, 1, $false, "", 0)
$findRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $findRange)

# --- Fill in the team member names table (3rd table in the document) ---
$table = $d.Tables.Item(3)
$table.Cell(2, 1).Range.Text = "Alex Norkus"
$table.Cell(3, 1).Range.Text = "Julius Mesa"
$table.Cell(4, 1).Range.Text = "Surya Manikhandan"
$table.Cell(5, 1).Range.Text = "Vincent Lin"
